# Zeitaufwand.xlsx - add a new log entry (row 5) and widen the "Tätigkeit" column.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 5 used to be a placeholder (only D5 had a style, no data). Fill it in:
#   B5 = 12.01.2018 (serial 43112), same date style as the rows above it
#   C5 = new shared-string description of the work done
#   D5 = 4 (hours)
# Copying B4's formatting first (instead of setting NumberFormat directly)
# makes the new cell reuse the existing date style rather than minting a
# duplicate one.
$ws.Range("B4").Copy()
$ws.Range("B5").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("B5").Value = 43112

$ws.Range("C5").Value = "Implementierung variabler Zitate auf Indexseite, Implementierung eines Videoarchivs, Modularisierung von Seiteninhalt. "

$ws.Range("D5").Value = 4

# D7's "SUM(D2:D6) & ""h""" formula will recalc on its own now that D5 has a
# value.

# Widen column C so the longer text fits.
$ws.Columns("C").ColumnWidth = 98.43

# Move the active-cell selection to C6, matching the saved view state.
$ws.Range("C6").Select()
